# Implementacion pruebas disparar Proyectil
$wb = $excel.ActiveWorkbook

# --- Sheet "Entregables": mark two deliverables as done (x) ---
$wsEntregables = $wb.Worksheets.Item("Entregables")
$wsEntregables.Range("B1").Value = "x"
$wsEntregables.Range("B5").Value = "x"
$wsEntregables.Range("B2").Select() | Out-Null

# --- Sheet "PruebasUnit": update projectile trajectory test data ---
$wsPruebas = $wb.Worksheets.Item("PruebasUnit")

$wsPruebas.Range("I1").Value = 50
$wsPruebas.Range("I2").Value = 50
$wsPruebas.Range("I4").Value = 100

$wsPruebas.Range("I6").Formula = "=I3-I1+(J3/2)"

$wsPruebas.Range("F7").Value = "x"
$wsPruebas.Range("F8").Value = "x"
$wsPruebas.Range("F9").Value = "x"

$wsPruebas.Range("I10").Formula = "=I18/I9"

$wsPruebas.Range("F30").Value = "x"
$wsPruebas.Range("F36").Value = "x"
$wsPruebas.Range("F37").Value = "x"

$wsPruebas.Range("D48").Select() | Out-Null

$wb.Save()
